$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 48177.047
$ws.Range("I28").Value = 53363.367
$ws.Range("J28").Value = 15330.333
$ws.Range("K28").Value = 53363.367
$ws.Range("L28").Value = 15330.333
$ws.Range("M28").Value = -52878.367
$ws.Range("N28").Value = -16300.333
$ws.Range("H86").Value = 22231100
$ws.Range("J86").Value = 40012480
$ws.Range("L86").Value = 40012480
$ws.Range("N86").Value = -40014726
$ws.Range("H89").Value = 22231100
$ws.Range("J89").Value = 40012480
$ws.Range("L89").Value = 200062400
$ws.Range("N89").Value = -200073632
$ws.Range("H100").Value = 13411.25
$ws.Range("I100").Value = 13411.25
$ws.Range("K100").Value = 13411.25
$ws.Range("M100").Value = -12870.25
$ws.Range("H106").Value = 118917.11
$ws.Range("I106").Value = 7792.3335
$ws.Range("K106").Value = 7792.3335
$ws.Range("M106").Value = -7161.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2993
$ws.Range("I32").Value = 2887.721
$ws.Range("K32").Value = 2887.721
$ws.Range("M32").Value = -2600.721
$ws.Range("H45").Value = 1971.2858
$ws.Range("I45").Value = 1559.8
$ws.Range("K45").Value = 1559.8
$ws.Range("M45").Value = -1182.8
$ws.Range("H63").Value = 200020050
$ws.Range("I63").Value = 333344300
$ws.Range("K63").Value = 333344300
$ws.Range("M63").Value = -333343614
$ws.Range("H66").Value = 200020050
$ws.Range("I66").Value = 333344300
$ws.Range("K66").Value = 1666721500
$ws.Range("M66").Value = -1666718068
$ws.Range("H102").Value = 6114995.5
$ws.Range("I102").Value = 10102566
$ws.Range("K102").Value = 10102566
$ws.Range("M102").Value = -10100944

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 55558276
$ws.Range("I80").Value = 3325
$ws.Range("K80").Value = 3325
$ws.Range("M80").Value = -2327
$ws.Range("H83").Value = 55558276
$ws.Range("I83").Value = 3325
$ws.Range("K83").Value = 16625
$ws.Range("M83").Value = -11633
$ws.Range("H94").Value = 2124.8667
$ws.Range("I94").Value = 1124.8182
$ws.Range("K94").Value = 1124.8182
$ws.Range("M94").Value = -673.8181999999999
$ws.Range("H134").Value = 6099718
$ws.Range("I134").Value = 6946429
$ws.Range("K134").Value = 20839287
$ws.Range("M134").Value = -20836752

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 93.85714
$ws.Range("I22").Value = 92.5
$ws.Range("K22").Value = 92.5
$ws.Range("M22").Value = 257.5
$ws.Range("H31").Value = 2946.726
$ws.Range("I31").Value = 2153.9119
$ws.Range("J31").Value = 3485.84
$ws.Range("K31").Value = 2153.9119
$ws.Range("L31").Value = 3485.84
$ws.Range("M31").Value = -1858.9119
$ws.Range("N31").Value = -4075.84
$ws.Range("H34").Value = 2946.726
$ws.Range("I34").Value = 2153.9119
$ws.Range("J34").Value = 3485.84
$ws.Range("K34").Value = 2153.9119
$ws.Range("L34").Value = 3485.84
$ws.Range("M34").Value = -1951.9119
$ws.Range("N34").Value = -3889.84
$ws.Range("H122").Value = 2656
$ws.Range("I122").Value = 2269.4707
$ws.Range("K122").Value = 6808.4121
$ws.Range("M122").Value = -4358.4121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3746.7646
$ws.Range("I62").Value = 3243.889
$ws.Range("K62").Value = 9731.667000000001
$ws.Range("M62").Value = -9045.667000000001
$ws.Range("H63").Value = 16989.666
$ws.Range("I63").Value = 16587.6
$ws.Range("J63").Value = 19000
$ws.Range("K63").Value = 49762.8
$ws.Range("L63").Value = 57000
$ws.Range("M63").Value = -49013.8
$ws.Range("N63").Value = -58498
$ws.Range("H64").Value = 4753.25
$ws.Range("I64").Value = 4999.6665
$ws.Range("J64").Value = 4014
$ws.Range("K64").Value = 14998.9995
$ws.Range("L64").Value = 12042
$ws.Range("M64").Value = -14728.9995
$ws.Range("N64").Value = -12582
$ws.Range("H65").Value = 3746.7646
$ws.Range("I65").Value = 3243.889
$ws.Range("K65").Value = 29195.001
$ws.Range("M65").Value = -25763.001
$ws.Range("H66").Value = 16989.666
$ws.Range("I66").Value = 16587.6
$ws.Range("J66").Value = 19000
$ws.Range("K66").Value = 149288.4
$ws.Range("L66").Value = 171000
$ws.Range("M66").Value = -145544.4
$ws.Range("N66").Value = -178488
$ws.Range("H67").Value = 4753.25
$ws.Range("I67").Value = 4999.6665
$ws.Range("J67").Value = 4014
$ws.Range("K67").Value = 14998.9995
$ws.Range("L67").Value = 12042
$ws.Range("M67").Value = -14062.9995
$ws.Range("N67").Value = -13914
$ws.Range("H75").Value = 1763
$ws.Range("I75").Value = 395
$ws.Range("J75").Value = 4499
$ws.Range("K75").Value = 1185
$ws.Range("L75").Value = 13497
$ws.Range("M75").Value = -187
$ws.Range("N75").Value = -15493
$ws.Range("H76").Value = 7149
$ws.Range("I76").Value = 3468.5
$ws.Range("K76").Value = 10405.5
$ws.Range("M76").Value = -10022.5
$ws.Range("H78").Value = 1763
$ws.Range("I78").Value = 395
$ws.Range("J78").Value = 4499
$ws.Range("K78").Value = 3555
$ws.Range("L78").Value = 40491
$ws.Range("M78").Value = 1437
$ws.Range("N78").Value = -50475
$ws.Range("H79").Value = 7149
$ws.Range("I79").Value = 3468.5
$ws.Range("K79").Value = 10405.5
$ws.Range("M79").Value = -9079.5
$ws.Range("H87").Value = 999.5
$ws.Range("I87").Value = 999.5
$ws.Range("K87").Value = 2998.5
$ws.Range("M87").Value = -1750.5
$ws.Range("H90").Value = 999.5
$ws.Range("I90").Value = 999.5
$ws.Range("K90").Value = 8995.5
$ws.Range("M90").Value = -2755.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 32500
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 32500
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 32500
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -32868
$ws.Range("H132").Value = 3647.44
$ws.Range("I132").Value = 2762.0625
$ws.Range("K132").Value = 8286.1875
$ws.Range("M132").Value = -5756.1875
$ws.Range("H139").Value = 212000
$ws.Range("J139").Value = 212000
$ws.Range("L139").Value = 212000
$ws.Range("N139").Value = -222280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1738.25
$ws.Range("I93").Value = 1738.25
$ws.Range("K93").Value = 1738.25
$ws.Range("M93").Value = -490.25
$ws.Range("H100").Value = 3413.1428
$ws.Range("I100").Value = 3078.4
$ws.Range("K100").Value = 3078.4
$ws.Range("M100").Value = -2537.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1504904.4
$ws.Range("J4").Value = 116551
$ws.Range("L4").Value = 116551
$ws.Range("N4").Value = -116777
$ws.Range("H5").Value = 5001.5
$ws.Range("I5").Value = 5001.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 5001.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4889.5
$ws.Range("N5").ClearContents()
$ws.Range("H81").Value = 7147851
$ws.Range("I81").Value = 2271.5293
$ws.Range("J81").Value = 18191018
$ws.Range("K81").Value = 4543.0586
$ws.Range("L81").Value = 36382036
$ws.Range("M81").Value = -3482.0586
$ws.Range("N81").Value = -36384158
$ws.Range("H84").Value = 7147851
$ws.Range("I84").Value = 2271.5293
$ws.Range("J84").Value = 18191018
$ws.Range("K84").Value = 22715.293
$ws.Range("L84").Value = 181910180
$ws.Range("M84").Value = -17411.293
$ws.Range("N84").Value = -181920788
$ws.Range("H100").Value = 1079.0834
$ws.Range("I100").Value = 935.5714
$ws.Range("K100").Value = 1871.1428
$ws.Range("M100").Value = -1330.1428
$ws.Range("H111").Value = 105000
$ws.Range("J111").Value = 105000
$ws.Range("L111").Value = 105000
$ws.Range("N111").Value = -113180
$ws.Range("H126").Value = 1608.7273
$ws.Range("I126").Value = 1519.6
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4558.799999999999
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2088.799999999999
$ws.Range("N126").Value = -12440
$ws.Range("H132").Value = 2344.818
$ws.Range("I132").Value = 2046.8
$ws.Range("K132").Value = 6140.4
$ws.Range("M132").Value = -3610.4
$ws.Range("H138").Value = 93997.5
$ws.Range("J138").Value = 93997.5
$ws.Range("L138").Value = 93997.5
$ws.Range("N138").Value = -104277.5
